$d = $word.ActiveDocument

$newText = 'Esteu participant en una campanya mundial per observar i anotar la brillantor de les estrelles més febles que es poden veure, com a mitjà per mesurar la contaminació lumínica en un lloc determinat. Localitzant i observant la  Constel·lació d''Hèrcules a la nit i comparant la brillantor de les estrelles del cel amb la brillantor que indiquen els mapes, gent de tot el món aprendran com els llums de la seva zona contribueixen a augmentar la contaminació lumínica. Les vostres aportacions a la base de dades activa faran palesa la visibilitat del cel nocturn.'

$found = $false
foreach ($p in $d.Paragraphs) {
    $t = $p.Range.Text
    if ($t -like "*Esteu participant*en*una*campanya*mundial*") {
        $r = $p.Range
        # Drop the paragraph mark from the range so we only touch the run content.
        $r.MoveEnd(1, -1) | Out-Null
        # Remove all the existing (separately formatted) runs in this paragraph...
        $r.Delete()
        # ...and insert one single, plain run with the merged/updated text
        # (no explicit run formatting, so it simply inherits the paragraph/style formatting).
        $ins = $p.Range
        $ins.Collapse(1) | Out-Null
        $ins.InsertAfter($newText)
        $found = $true
        break
    }
}

if (-not $found) {
    Write-Host "ERROR: target paragraph not found"
} else {
    Write-Host "Paragraph replaced successfully"
}
